# Populate the "Responsibile Parties" sheet (CPTEC-INPE responsible-party
# records) and make it the active sheet, matching the authored workbook.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Responsibile Parties")

# --- Row heights for the three new data rows -------------------------------
$ws3.Rows.Item(3).RowHeight = 42
$ws3.Rows.Item(4).RowHeight = 40
$ws3.Rows.Item(5).RowHeight = 30

# --- Row 3: CPTEC-INPE (organisation record) --------------------------------
$ws3.Range("A3").Value = "CPTEC-INPE"
$ws3.Range("B3").Value = "Center for weather forecast and Climate studies - National Institute for Space Research"
$ws3.Range("C3").Value = $true
$ws3.Range("D3").Value = "Rod. Pres. Dutra km39, Cahcoeira Paulista, SP - Brasil - CEP 12630-000"
$ws3.Range("E3").Value = "besm@inpe.br"
$ws3.Hyperlinks.Add($ws3.Range("E3"), "mailto:besm@inpe.br") | Out-Null
$ws3.Range("F3").Value = "http://www.inpe.br/besm/"

# --- Row 4: André Lanfer (individual record) --------------------------------
$ws3.Range("A4").Value = "André-Lanfer"
$ws3.Range("B4").Value = "André Lanfer"
$ws3.Range("C4").Value = $false
$ws3.Range("D4").Value = "Rod. Pres. Dutra km39, Cahcoeira Paulista, SP - Brasil - CEP 12630-001"
$ws3.Range("D4").Borders.Item(8).LineStyle = 0
$ws3.Range("D4").HorizontalAlignment = -4131
$ws3.Range("E4").Value = "andre.lanfer@inpe.br"
$ws3.Hyperlinks.Add($ws3.Range("E4"), "mailto:andre.lanfer@inpe.br") | Out-Null
$ws3.Range("F4").Value = "http://www.inpe.br/besm/equipe/"
$ws3.Range("G4").Value = "https://orcid.org/0000-0001-8719-8045"

# --- Row 5: Paulo Nobre (individual record) ---------------------------------
$ws3.Range("A5").Value = "Paulo Nobre"
$ws3.Range("B5").Value = "Paulo Nobre"
$ws3.Range("C5").Value = $false
$ws3.Range("D5").Value = "Rod. Pres. Dutra km39, Cahcoeira Paulista, SP - Brasil - CEP 12630-002"
$ws3.Range("D5").Borders.Item(8).LineStyle = 0
$ws3.Range("E5").Value = "paulo.nobre@inpe.br"
$ws3.Hyperlinks.Add($ws3.Range("E5"), "mailto:paulo.nobre@inpe.br") | Out-Null
$ws3.Range("F5").Value = "http://www.inpe.br/besm/equipe/"
$ws3.Range("F5").HorizontalAlignment = 1

# --- Selection / active sheet -----------------------------------------------
$ws3.Activate()
$ws3.Range("F4:F5").Select()
